{"js": "// Collapse the multiple single-word/space runs in the Title, Author and\n// Abstract paragraphs into one run per paragraph containing the full text.\n// The visible text is unchanged -- only the run structure is simplified.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"style\");\nawait context.sync();\n\nconst targets = {\n  \"Title\": \"Questions: Introduction to sigma notation\",\n  \"Author\": \"Ifan Howells-Baines, Mark Toner\",\n  \"Abstract\": \"Questions relating to the guide on introduction to sigma notation.\"\n};\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const newText = targets[paragraph.style];\n  if (newText !== undefined) {\n    paragraph.getRange().insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Collapse the multiple single-word/space runs in the Title, Author and\n# Abstract paragraphs into one run per paragraph containing the full text.\n# The visible text is unchanged -- only the run structure is simplified.\n\n$d = $word.ActiveDocument\n\n$targets = @{\n    \"Title\"    = \"Questions: Introduction to sigma notation\"\n    \"Author\"   = \"Ifan Howells-Baines, Mark Toner\"\n    \"Abstract\" = \"Questions relating to the guide on introduction to sigma notation.\"\n}\n\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Style.NameLocal\n    if ($targets.ContainsKey($styleName)) {\n        $newText = $targets[$styleName]\n\n        $r = $p.Range\n        [void]$r.MoveEnd(1, -1)   # exclude the trailing paragraph mark\n\n        $find = $r.Find\n        $find.ClearFormatting()\n        $find.Replacement.ClearFormatting()\n        $find.Text = $r.Text\n        $find.Replacement.Text = $newText\n        [void]$find.Execute($r.Text, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    }\n}\n"}
